$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 14.817691
$ws.Range("N2").Value = 29.635382
$ws.Range("O2").Value = 0.2230198506330352
$ws.Range("P2").Value = 0.1687294465231367
$ws.Range("Q2").Value = 0.4146039334103333
$ws.Range("R2").Value = 2.487623600462
$ws.Range("S2").Value = 0.2230198506330352
$ws.Range("T2").Value = 0.1687294465231367

# Row 3
$ws.Range("O3").Value = 0.3961505993138794
$ws.Range("P3").Value = 0.4495716715730233
$ws.Range("S3").Value = 0.3961505993138794
$ws.Range("T3").Value = 0.4495716715730233

# Row 4
$ws.Range("M4").Value = 8.453176999999998
$ws.Range("N4").Value = 25.359531
$ws.Range("O4").Value = 0.1272280729780779
$ws.Range("P4").Value = 0.144384831270821
$ws.Range("Q4").Value = 0.2365227101856666
$ws.Range("R4").Value = 2.128704391671
$ws.Range("S4").Value = 0.1272280729780779
$ws.Range("T4").Value = 0.144384831270821

# Row 5
$ws.Range("M5").Value = 8.867229500000001
$ws.Range("N5").Value = 17.734459
$ws.Range("O5").Value = 0.133459943159757
$ws.Range("P5").Value = 0.1009713811503176
$ws.Range("Q5").Value = 0.2481080371531667
$ws.Range("R5").Value = 1.488648222919
$ws.Range("S5").Value = 0.133459943159757
$ws.Range("T5").Value = 0.1009713811503176

# Row 6
$ws.Range("M6").Value = 3.775572333333333
$ws.Range("N6").Value = 11.326717
$ws.Range("O6").Value = 0.05682582919526532
$ws.Range("P6").Value = 0.06448881577886201
$ws.Range("Q6").Value = 0.1056417724107778
$ws.Range("R6").Value = 0.9507759516969999
$ws.Range("S6").Value = 0.05682582919526532
$ws.Range("T6").Value = 0.06448881577886201

# Row 7
$ws.Range("M7").Value = 4.206767
$ws.Range("N7").Value = 12.620301
$ws.Range("O7").Value = 0.06331570471998517
$ws.Range("P7").Value = 0.07185385370383916
$ws.Range("Q7").Value = 0.1177067429156667
$ws.Range("R7").Value = 1.059360686241
$ws.Range("S7").Value = 0.06331570471998517
$ws.Range("T7").Value = 0.07185385370383916
